$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 186262.312244738
$ws.Range("C2").Value = 270117.1320619749
$ws.Range("D2").Value = 316273.2183627473
$ws.Range("E2").Value = 341341.5844877228
$ws.Range("B3").Value = 225117.3819338269
$ws.Range("C3").Value = 323747.5607002651
$ws.Range("D3").Value = 373956.5822425464
$ws.Range("E3").Value = 405779.931373395
$ws.Range("B4").Value = 201010.1057277275
$ws.Range("C4").Value = 299609.6048645193
$ws.Range("D4").Value = 354511.2290197113
$ws.Range("E4").Value = 387666.0245078456
$ws.Range("B5").Value = 155394.6062134071
$ws.Range("C5").Value = 219786.4634951173
$ws.Range("D5").Value = 246916.6157535585
$ws.Range("E5").Value = 269805.9494044766
$ws.Range("B6").Value = 136716.9600814675
$ws.Range("C6").Value = 192827.5739427786
$ws.Range("D6").Value = 218716.889540656
$ws.Range("E6").Value = 236080.4277308756
$ws.Range("B7").Value = 14675.8668519096
$ws.Range("C7").Value = 20468.52016431294
$ws.Range("D7").Value = 23275.2447736936
$ws.Range("E7").Value = 24997.44314802386
$ws.Range("B8").Value = 705765.2691399819
$ws.Range("C8").Value = 1030473.0229735
$ws.Range("D8").Value = 1228209.609252516
$ws.Range("E8").Value = 1334096.69295451
$ws.Range("B9").Value = 199454.2841605958
$ws.Range("C9").Value = 280426.2657429983
$ws.Range("D9").Value = 318482.8156063193
$ws.Range("E9").Value = 344244.995487882
$ws.Range("B10").Value = 85634.35353239594
$ws.Range("C10").Value = 117364.6452962163
$ws.Range("D10").Value = 134775.7304153467
$ws.Range("E10").Value = 142995.7798438304
$ws.Range("B11").Value = 15740.31468564153
$ws.Range("C11").Value = 20485.15811294451
$ws.Range("D11").Value = 23259.24230106769
$ws.Range("E11").Value = 26626.23532188547
$ws.Range("B12").Value = 37527.28018678263
$ws.Range("C12").Value = 53768.00001602747
$ws.Range("D12").Value = 61460.37232183362
$ws.Range("E12").Value = 64361.18566302748
$ws.Range("B13").Value = 48128.78275800219
$ws.Range("C13").Value = 65393.02664912661
$ws.Range("D13").Value = 75816.81769483958
$ws.Range("E13").Value = 81388.17461996997
